$d = $word.ActiveDocument

$pairs = @(
    @("434×9=", "848×7="),
    @("714×5=", "857×5="),
    @("439×3=", "684×5="),
    @("449×7=", "665×7="),
    @("673×6=", "176×5="),
    @("435×2=", "796×2="),
    @("603×4=", "592×7="),
    @("385×4=", "661×9="),
    @("915×7=", "972×9="),
    @("495×4=", "657×4="),
    @("263×4=", "191×8="),
    @("510×5=", "524×2="),
    @("758×4=", "405×2="),
    @("267×4=", "740×2="),
    @("129×4=", "611×9="),
    @("266×4=", "652×3="),
    @("572×2=", "394×9="),
    @("179×2=", "853×3="),
    @("434×5=", "641×2="),
    @("765×8=", "475×5="),
    @("197×8=", "458×9="),
    @("835×5=", "717×6="),
    @("195×4=", "372×2="),
    @("553×2=", "479×6="),
    @("479×4=", "152×7=")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
